$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "Datos actualizados" timestamp string (row 1, col A)
$ws.Range("A1").Value = "Datos actualizados a 17 de Septiembre de 2020 a las 08:33"

# Row 28 - Ucrania
$ws.Range("B28").Value = 166244
$ws.Range("C28").Value = 3584
$ws.Range("D28").Value = 73913
$ws.Range("E28").Value = 88931
$ws.Range("G28").Value = 60
$ws.Range("H28").Value = 3400

# Row 59 - Uzbekistan
$ws.Range("B59").Value = 49385
$ws.Range("C59").Value = 370
$ws.Range("E59").Value = 3500
$ws.Range("G59").Value = 4
$ws.Range("H59").Value = 411

# Row 64 - Kirguistan
$ws.Range("B64").Value = 45153
$ws.Range("C64").Value = 81
$ws.Range("D64").Value = 41317
$ws.Range("E64").Value = 2773

# Row 75 - El Salvador
$ws.Range("E75").Value = 6402
$ws.Range("G75").Value = 5
$ws.Range("H75").Value = 801
